# Apply row-content rotation per commit diff: rows 4,5,10,11,15,16,19-21,
# 24-26,27-30,31-34,35-38,39,40,43,44,48-50 each take on the content that
# used to live in another row of the same cycle (full record swap).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowData = @{}

$rowData[4] = @(130754287, 57881, $null, 'NT', 100049, 'Spillkråka', 'Dryocopus martius', '(Linnaeus, 1758)', $null, $null, $null, $null, 'äldre spår', $null, $null, 'Kråkbackarna, Dlr', 490501, 6763773, 10, 'Dalarna', 'Mora', 'Dalarna', 'Mora', $null, $null, '11:43', $null, '11:43', $null, $false, $false, $null, $false, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 'Håkan Thenander', 'Håkan Thenander, Bo karlstens', $null)
$rowData[5] = @(130789471, 79243, $null, 'NT', 6425, 'Garnlav', 'Alectoria sarmentosa', '(Ach.) Ach.', $null, $null, $null, $null, $null, $null, $null, 'Kråkbackarna, Dlr', 490498, 6763669, 10, 'Dalarna', 'Mora', 'Dalarna', 'Mora', $null, $null, '14:18', $null, '14:18', $null, $false, $false, $null, $false, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 'Bo karlstens', 'Bo karlstens, Håkan Thenander', $null)
$rowData[10] = @(130789482, 79243, $null, 'NT', 6425, 'Garnlav', 'Alectoria sarmentosa', '(Ach.) Ach.', $null, $null, $null, $null, $null, $null, $null, 'Kråkbackarna, Dlr', 490451, 6764011, 10, 'Dalarna', 'Mora', 'Dalarna', 'Mora', $null, $null, '13:53', $null, '13:53', $null, $false, $false, $null, $false, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 'Bo karlstens', 'Bo karlstens, Håkan Thenander', $null)
$rowData[11] = @(130789477, 79243, $null, 'NT', 6425, 'Garnlav', 'Alectoria sarmentosa', '(Ach.) Ach.', $null, $null, $null, $null, $null, $null, $null, 'Kråkbackarna, Dlr', 490435, 6764047, 10, 'Dalarna', 'Mora', 'Dalarna', 'Mora', $null, $null, '14:00', $null, '14:00', $null, $false, $false, $null, $false, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 'Bo karlstens', 'Bo karlstens, Håkan Thenander', $null)
$rowData[15] = @(130754083, 79243, $null, 'NT', 6425, 'Garnlav', 'Alectoria sarmentosa', '(Ach.) Ach.', $null, $null, $null, $null, $null, $null, $null, 'Kråkbackarna, Dlr', 490529, 6763665, 10, 'Dalarna', 'Mora', 'Dalarna', 'Mora', $null, $null, '11:43', $null, '11:43', '1 bild, tall med gran till vänster', $false, $false, $null, $false, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 'Håkan Thenander', 'Håkan Thenander, Bo karlstens', $null)
$rowData[16] = @(130789470, 79243, $null, 'NT', 6425, 'Garnlav', 'Alectoria sarmentosa', '(Ach.) Ach.', $null, $null, $null, $null, $null, $null, $null, 'Kråkbackarna, Dlr', 490376, 6763596, 10, 'Dalarna', 'Mora', 'Dalarna', 'Mora', $null, $null, '15:10', $null, '15:10', $null, $false, $false, $null, $false, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 'Bo karlstens', 'Bo karlstens, Håkan Thenander', $null)
$rowData[19] = @(130757199, 79243, $null, 'NT', 6425, 'Garnlav', 'Alectoria sarmentosa', '(Ach.) Ach.', $null, $null, $null, $null, $null, $null, $null, 'Kråkbackarna, Dlr', 490467, 6763573, 10, 'Dalarna', 'Mora', 'Dalarna', 'Mora', $null, $null, '15:01', $null, '15:01', $null, $false, $false, $null, $false, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 'Håkan Thenander', 'Håkan Thenander, Bo karlstens', $null)
$rowData[20] = @(130757425, 79243, $null, 'NT', 6425, 'Garnlav', 'Alectoria sarmentosa', '(Ach.) Ach.', $null, $null, $null, $null, $null, $null, $null, 'Kråkbackarna, Dlr', 490361, 6763574, 10, 'Dalarna', 'Mora', 'Dalarna', 'Mora', $null, $null, '15:01', $null, '15:01', $null, $false, $false, $null, $false, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 'Håkan Thenander', 'Håkan Thenander, Bo karlstens', $null)
$rowData[21] = @(130754307, 79243, $null, 'NT', 6425, 'Garnlav', 'Alectoria sarmentosa', '(Ach.) Ach.', $null, $null, $null, $null, $null, $null, $null, 'Kråkbackarna, Dlr', 490501, 6763773, 10, 'Dalarna', 'Mora', 'Dalarna', 'Mora', $null, $null, '11:43', $null, '11:43', $null, $false, $false, $null, $false, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 'Håkan Thenander', 'Håkan Thenander, Bo karlstens', $null)
$rowData[24] = @(130754985, 79243, $null, 'NT', 6425, 'Garnlav', 'Alectoria sarmentosa', '(Ach.) Ach.', $null, $null, $null, $null, $null, $null, $null, 'Prikattmyren, Dlr', 490434, 6764045, 10, 'Dalarna', 'Mora', 'Dalarna', 'Mora', $null, $null, '11:43', $null, '11:43', $null, $false, $false, $null, $false, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 'Håkan Thenander', 'Håkan Thenander, Bo karlstens', $null)
$rowData[25] = @(130754896, 79243, $null, 'NT', 6425, 'Garnlav', 'Alectoria sarmentosa', '(Ach.) Ach.', $null, $null, $null, $null, $null, $null, $null, 'Prikattmyren, Dlr', 490446, 6764008, 10, 'Dalarna', 'Mora', 'Dalarna', 'Mora', $null, $null, '11:43', $null, '11:43', $null, $false, $false, $null, $false, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 'Håkan Thenander', 'Håkan Thenander, Bo karlstens', $null)
$rowData[26] = @(130789495, 79243, $null, 'NT', 6425, 'Garnlav', 'Alectoria sarmentosa', '(Ach.) Ach.', $null, $null, $null, $null, $null, $null, $null, 'Kråkbackarna, Dlr', 490496, 6763849, 10, 'Dalarna', 'Mora', 'Dalarna', 'Mora', $null, $null, '13:33', $null, '13:33', $null, $false, $false, $null, $false, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 'Bo karlstens', 'Bo karlstens, Håkan Thenander', $null)
$rowData[27] = @(130757159, 79243, $null, 'NT', 6425, 'Garnlav', 'Alectoria sarmentosa', '(Ach.) Ach.', $null, $null, $null, $null, $null, $null, $null, 'Kråkbackarna, Dlr', 490482, 6763574, 10, 'Dalarna', 'Mora', 'Dalarna', 'Mora', $null, $null, '15:01', $null, '15:01', $null, $false, $false, $null, $false, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 'Håkan Thenander', 'Håkan Thenander, Bo karlstens', $null)
$rowData[28] = @(130789472, 79243, $null, 'NT', 6425, 'Garnlav', 'Alectoria sarmentosa', '(Ach.) Ach.', $null, $null, $null, $null, $null, $null, $null, 'Kråkbackarna, Dlr', 490450, 6763926, 10, 'Dalarna', 'Mora', 'Dalarna', 'Mora', $null, $null, '14:07', $null, '14:07', $null, $false, $false, $null, $false, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 'Bo karlstens', 'Bo karlstens, Håkan Thenander', $null)
$rowData[29] = @(130754796, 79243, $null, 'NT', 6425, 'Garnlav', 'Alectoria sarmentosa', '(Ach.) Ach.', $null, $null, $null, $null, $null, $null, $null, 'Kråkbackarna, Dlr', 490463, 6763939, 10, 'Dalarna', 'Mora', 'Dalarna', 'Mora', $null, $null, '11:43', $null, '11:43', $null, $false, $false, $null, $false, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 'Håkan Thenander', 'Håkan Thenander, Bo karlstens', $null)
$rowData[30] = @(130755667, 57884, $null, 'NT', 100109, 'Tretåig hackspett', 'Picoides tridactylus', '(Linnaeus, 1758)', $null, $null, $null, $null, 'färska spår', $null, $null, 'Prikattmyren, Dlr', 490444, 6763770, 10, 'Dalarna', 'Mora', 'Dalarna', 'Mora', $null, $null, '11:43', $null, '11:43', '2 bilder', $false, $false, $null, $false, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 'Håkan Thenander', 'Håkan Thenander, Bo karlstens', $null)
$rowData[31] = @(130789475, 79243, $null, 'NT', 6425, 'Garnlav', 'Alectoria sarmentosa', '(Ach.) Ach.', $null, $null, $null, $null, $null, $null, $null, 'Kråkbackarna, Dlr', 490441, 6763984, 10, 'Dalarna', 'Mora', 'Dalarna', 'Mora', $null, $null, '14:03', $null, '14:03', $null, $false, $false, $null, $false, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 'Bo karlstens', 'Bo karlstens, Håkan Thenander', $null)
$rowData[32] = @(130757236, 5177, $null, 'LC', 100526, 'Bronshjon', 'Callidium coriaceum', 'Paykull, 1800', $null, $null, $null, $null, 'äldre gnagspår', $null, $null, 'Kråkbackarna, Dlr', 490467, 6763573, 10, 'Dalarna', 'Mora', 'Dalarna', 'Mora', $null, $null, '15:01', $null, '15:01', $null, $false, $false, $null, $false, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 'Håkan Thenander', 'Håkan Thenander, Bo karlstens', $null)
$rowData[33] = @(130757715, 8451, $null, 'LC', 106545, 'Mindre märgborre', 'Tomicus minor', '(Hartig, 1834)', $null, $null, $null, $null, 'äldre gnagspår', $null, $null, 'Brunnvasselänget, Dlr', 490186, 6763602, 10, 'Dalarna', 'Mora', 'Dalarna', 'Mora', $null, $null, '15:01', $null, '15:01', $null, $false, $false, $null, $false, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 'Håkan Thenander', 'Håkan Thenander, Bo karlstens', $null)
$rowData[34] = @(130789476, 79243, $null, 'NT', 6425, 'Garnlav', 'Alectoria sarmentosa', '(Ach.) Ach.', $null, $null, $null, $null, $null, $null, $null, 'Kråkbackarna, Dlr', 490438, 6764018, 10, 'Dalarna', 'Mora', 'Dalarna', 'Mora', $null, $null, '14:02', $null, '14:02', $null, $false, $false, $null, $false, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 'Bo karlstens', 'Bo karlstens, Håkan Thenander', $null)
$rowData[35] = @(130757247, 5197, $null, 'LC', 105930, 'Vågbandad barkbock', 'Semanotus undatus', '(Linnaeus, 1758)', $null, $null, $null, $null, 'äldre gnagspår', $null, $null, 'Kråkbackarna, Dlr', 490467, 6763573, 10, 'Dalarna', 'Mora', 'Dalarna', 'Mora', $null, $null, '15:01', $null, '15:01', $null, $false, $false, $null, $false, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 'Håkan Thenander', 'Håkan Thenander, Bo karlstens', $null)
$rowData[36] = @(130754851, 79243, $null, 'NT', 6425, 'Garnlav', 'Alectoria sarmentosa', '(Ach.) Ach.', $null, $null, $null, $null, $null, $null, $null, 'Kråkbackarna, Dlr', 490449, 6763949, 10, 'Dalarna', 'Mora', 'Dalarna', 'Mora', $null, $null, '11:43', $null, '11:43', '1 bild. På tallstam', $false, $false, $null, $false, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 'Håkan Thenander', 'Håkan Thenander, Bo karlstens', $null)
$rowData[37] = @(130757412, 79243, $null, 'NT', 6425, 'Garnlav', 'Alectoria sarmentosa', '(Ach.) Ach.', $null, $null, $null, $null, $null, $null, $null, 'Kråkbackarna, Dlr', 490381, 6763583, 10, 'Dalarna', 'Mora', 'Dalarna', 'Mora', $null, $null, '15:01', $null, '15:01', $null, $false, $false, $null, $false, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 'Håkan Thenander', 'Håkan Thenander, Bo karlstens', $null)
$rowData[38] = @(130789468, 79243, $null, 'NT', 6425, 'Garnlav', 'Alectoria sarmentosa', '(Ach.) Ach.', $null, $null, $null, $null, $null, $null, $null, 'Kråkbackarna, Dlr', 490321, 6763593, 10, 'Dalarna', 'Mora', 'Dalarna', 'Mora', $null, $null, '15:14', $null, '15:14', $null, $false, $false, $null, $false, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 'Bo karlstens', 'Bo karlstens, Håkan Thenander', $null)
$rowData[39] = @(130789489, 79243, $null, 'NT', 6425, 'Garnlav', 'Alectoria sarmentosa', '(Ach.) Ach.', $null, $null, $null, $null, $null, $null, $null, 'Kråkbackarna, Dlr', 490467, 6763913, 10, 'Dalarna', 'Mora', 'Dalarna', 'Mora', $null, $null, '13:40', $null, '13:40', $null, $false, $false, $null, $false, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 'Bo karlstens', 'Bo karlstens, Håkan Thenander', $null)
$rowData[40] = @(130789487, 79243, $null, 'NT', 6425, 'Garnlav', 'Alectoria sarmentosa', '(Ach.) Ach.', $null, $null, $null, $null, $null, $null, $null, 'Kråkbackarna, Dlr', 490472, 6763924, 10, 'Dalarna', 'Mora', 'Dalarna', 'Mora', $null, $null, '13:42', $null, '13:42', $null, $false, $false, $null, $false, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 'Bo karlstens', 'Bo karlstens, Håkan Thenander', $null)
$rowData[43] = @(130755062, 57884, $null, 'NT', 100109, 'Tretåig hackspett', 'Picoides tridactylus', '(Linnaeus, 1758)', $null, $null, $null, $null, 'äldre spår', $null, $null, 'Prikattmyren, Dlr', 490434, 6764045, 10, 'Dalarna', 'Mora', 'Dalarna', 'Mora', $null, $null, '11:43', $null, '11:43', '1 bild', $false, $false, $null, $false, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 'Håkan Thenander', 'Håkan Thenander, Bo karlstens', $null)
$rowData[44] = @(130757564, 79243, $null, 'NT', 6425, 'Garnlav', 'Alectoria sarmentosa', '(Ach.) Ach.', $null, $null, $null, $null, $null, $null, $null, 'Brunnvasselänget, Dlr', 490230, 6763582, 10, 'Dalarna', 'Mora', 'Dalarna', 'Mora', $null, $null, '15:01', $null, '15:01', $null, $false, $false, $null, $false, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 'Håkan Thenander', 'Håkan Thenander, Bo karlstens', $null)
$rowData[48] = @(130789459, 79243, $null, 'NT', 6425, 'Garnlav', 'Alectoria sarmentosa', '(Ach.) Ach.', $null, $null, $null, $null, $null, $null, $null, 'Kråkbackarna, Dlr', 490499, 6763662, 10, 'Dalarna', 'Mora', 'Dalarna', 'Mora', $null, $null, '14:27', $null, '14:27', 'vid gran 250 år plus', $false, $false, $null, $false, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 'Bo karlstens', 'Bo karlstens, Håkan Thenander', $null)
$rowData[49] = @(130758032, 79243, $null, 'NT', 6425, 'Garnlav', 'Alectoria sarmentosa', '(Ach.) Ach.', $null, $null, $null, $null, $null, $null, $null, 'Brunnvasselänget, Dlr', 490175, 6763613, 10, 'Dalarna', 'Mora', 'Dalarna', 'Mora', $null, $null, '15:01', $null, '15:01', $null, $false, $false, $null, $false, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 'Håkan Thenander', 'Håkan Thenander, Bo karlstens', $null)
$rowData[50] = @(130758572, 57881, $null, 'NT', 100049, 'Spillkråka', 'Dryocopus martius', '(Linnaeus, 1758)', $null, $null, $null, $null, 'färska spår', $null, $null, 'Kråkbackarna, Dlr', 490494, 6763540, 10, 'Dalarna', 'Mora', 'Dalarna', 'Mora', $null, $null, '15:01', $null, '15:01', $null, $false, $false, $null, $false, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 'Håkan Thenander', 'Håkan Thenander, Bo karlstens', $null)

foreach ($r in $rowData.Keys) {
    $vals = $rowData[$r]
    $arr = New-Object "object[,]" 1,51
    for ($i = 0; $i -lt 51; $i++) {
        $arr[0,$i] = $vals[$i]
    }
    $ws.Range("A" + $r + ":AY" + $r).Value2 = $arr
}

# Restore Y/AA (Startdatum/Slutdatum) as literal text "2026-01-18" without
# letting Excel COM reinterpret the ISO-looking string as a real date.
$dateRows = @(4, 5, 10, 11, 15, 16, 19, 20, 21, 24, 25, 26, 27, 28, 29, 30, 31, 32, 33, 34, 35, 36, 37, 38, 39, 40, 43, 44, 48, 49, 50)
foreach ($r in $dateRows) {
    $cellY = $ws.Range("Y" + $r)
    $cellY.NumberFormat = "@"
    $cellY.Value = "2026-01-18"
    $cellY.NumberFormat = "General"
    $cellAA = $ws.Range("AA" + $r)
    $cellAA.NumberFormat = "@"
    $cellAA.Value = "2026-01-18"
    $cellAA.NumberFormat = "General"
}
